$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.896.78'
$ws.Range("E2").Value = '  +0.16%  '
$ws.Range("D3").Value = '2.295.83'
$ws.Range("E3").Value = '  -0.81%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = "'108.01"
$ws.Range("E5").Value = '  +11.06%  '
$ws.Range("D6").Value = "'271.20"
$ws.Range("E6").Value = '  -0.49%  '
$ws.Range("E7").Value = '  -1.67%  '
$ws.Range("E8").Value = '  +0.19%  '
$ws.Range("D9").Value = "'0.615"
$ws.Range("E9").Value = '  -1.80%  '
$ws.Range("D10").Value = "'47.05"
$ws.Range("E10").Value = '  +3.61%  '
$ws.Range("D11").Value = "'0.0936"
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("E12").Value = '  +4.34%  '
$ws.Range("E13").Value = '  +1.68%  '
$ws.Range("D14").Value = "'15.66"
$ws.Range("E14").Value = '  +1.04%  '
$ws.Range("D15").Value = '2.637.69'
$ws.Range("E15").Value = '  -0.43%  '
$ws.Range("D16").Value = "'0.857"
$ws.Range("E16").Value = '  -1.87%  '
$ws.Range("D17").Value = '2.294.97'
$ws.Range("E17").Value = '  -0.99%  '
$ws.Range("D18").Value = '43.794.25'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("E19").Value = '  +0.27%  '
$ws.Range("D20").Value = "'6.30"
$ws.Range("E20").Value = '  -2.00%  '
$ws.Range("D21").Value = "'72.19"
$ws.Range("E21").Value = '  -1.82%  '
$ws.Range("E22").Value = '  +8.91%  '
$ws.Range("D23").Value = "'233.67"
$ws.Range("E23").Value = '  -2.49%  '
$ws.Range("E24").Value = '  +15.75%  '
$ws.Range("D25").Value = "'9.32"
$ws.Range("E25").Value = '  -1.60%  '
$ws.Range("E26").Value = '  +0.07%  '
$ws.Range("D27").Value = "'11.32"
$ws.Range("E27").Value = '  -0.36%  '
$ws.Range("D28").Value = "'40.74"
$ws.Range("E28").Value = '  +6.41%  '
$ws.Range("E29").Value = '  -1.28%  '
$ws.Range("E30").Value = '  -0.42%  '
$ws.Range("D31").Value = "'177.99"
$ws.Range("E31").Value = '  +1.62%  '
$ws.Range("D32").Value = "'21.88"
$ws.Range("E32").Value = '  -2.51%  '
$ws.Range("D33").Value = "'0.0910"
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("D34").Value = "'5.56"
$ws.Range("E34").Value = '  +1.26%  '
$ws.Range("D35").Value = "'4.90"
$ws.Range("E35").Value = '  +9.81%  '
$ws.Range("E36").Value = '  -0.25%  '
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("D38").Value = "'0.0360"
$ws.Range("E38").Value = '  -1.55%  '
$ws.Range("D39").Value = "'3.64"
$ws.Range("E39").Value = '  +8.21%  '
$ws.Range("E40").Value = '  -3.38%  '
$ws.Range("D41").Value = "'2.33"
$ws.Range("E41").Value = '  -2.95%  '
$ws.Range("E42").Value = '  -3.76%  '
$ws.Range("D43").Value = "'66.34"
$ws.Range("E43").Value = '  +5.44%  '
$ws.Range("D44").Value = "'12.19"
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("D45").Value = "'5.51"
$ws.Range("E45").Value = '  +3.08%  '
$ws.Range("D46").Value = "'8.79"
$ws.Range("E46").Value = '  -4.34%  '
$ws.Range("D47").Value = "'0.102"
$ws.Range("E47").Value = '  -1.81%  '
$ws.Range("E48").Value = '  +2.25%  '
$ws.Range("D49").Value = "'99.35"
$ws.Range("E49").Value = '  -1.17%  '
$ws.Range("E50").Value = '  +11.19%  '
$ws.Range("D51").Value = "'0.441"
$ws.Range("E51").Value = '  +5.15%  '
